$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45919
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 36.348
$ws.Range("A3").Value = 45919.01041666666
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 20.838
$ws.Range("A4").Value = 45919.02083333334
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 4.661
$ws.Range("A5").Value = 45919.03125
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 4.753
$ws.Range("A6").Value = 45919.04166666666
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 36.567
$ws.Range("A7").Value = 45919.05208333334
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 4.865
$ws.Range("A8").Value = 45919.0625
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 4.085
$ws.Range("A9").Value = 45919.07291666666
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 2.98
$ws.Range("A10").Value = 45919.08333333334
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 4.714
$ws.Range("A11").Value = 45919.09375
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 2.165
$ws.Range("A12").Value = 45919.10416666666
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1.711
$ws.Range("A13").Value = 45919.11458333334
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0.8129999999999999
$ws.Range("A14").Value = 45919.125
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.216
$ws.Range("A15").Value = 45919.13541666666
$ws.Range("B15").Value = 0.008
$ws.Range("C15").Value = 0.437
$ws.Range("A16").Value = 45919.14583333334
$ws.Range("B16").Value = 1.622
$ws.Range("C16").Value = 2.447
$ws.Range("A17").Value = 45919.15625
$ws.Range("B17").Value = 2.223
$ws.Range("C17").Value = 1.975
$ws.Range("A18").Value = 45919.16666666666
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 1.592
$ws.Range("A19").Value = 45919.17708333334
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0.904
$ws.Range("A20").Value = 45919.1875
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 5.011
$ws.Range("A21").Value = 45919.19791666666
$ws.Range("B21").Value = 0.018
$ws.Range("C21").Value = 8.311
$ws.Range("A22").Value = 45919.20833333334
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 3.014
$ws.Range("A23").Value = 45919.21875
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 3.25
$ws.Range("A24").Value = 45919.22916666666
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 13.804
$ws.Range("A25").Value = 45919.23958333334
$ws.Range("B25").Value = 1.149
$ws.Range("C25").Value = 4.554
$ws.Range("A26").Value = 45919.25
$ws.Range("B26").Value = 0.006
$ws.Range("C26").Value = 6.841
$ws.Range("A27").Value = 45919.26041666666
$ws.Range("B27").Value = 0.001
$ws.Range("C27").Value = 5.123
$ws.Range("A28").Value = 45919.27083333334
$ws.Range("B28").Value = 0.001
$ws.Range("C28").Value = 14.975
$ws.Range("A29").Value = 45919.28125
$ws.Range("B29").Value = 1.225
$ws.Range("C29").Value = 10.858
$ws.Range("A30").Value = 45919.29166666666
$ws.Range("B30").Value = 12.994
$ws.Range("C30").Value = 0
$ws.Range("A31").Value = 45919.30208333334
$ws.Range("B31").Value = 23.054
$ws.Range("C31").Value = 0
$ws.Range("A32").Value = 45919.3125
$ws.Range("B32").Value = 15.669
$ws.Range("C32").Value = 0
$ws.Range("A33").Value = 45919.32291666666
$ws.Range("B33").Value = 14.885
$ws.Range("C33").Value = 0
$ws.Range("A34").Value = 45919.33333333334
$ws.Range("B34").Value = 19.995
$ws.Range("C34").Value = 0
$ws.Range("A35").Value = 45919.34375
$ws.Range("B35").Value = 52.76
$ws.Range("C35").Value = 0
$ws.Range("A36").Value = 45919.35416666666
$ws.Range("B36").Value = 56.14
$ws.Range("C36").Value = 0
$ws.Range("A37").Value = 45919.36458333334
$ws.Range("B37").Value = 26.283
$ws.Range("C37").Value = 0
$ws.Range("A38").Value = 45919.375
$ws.Range("B38").Value = 22.169
$ws.Range("C38").Value = 0
$ws.Range("A39").Value = 45919.38541666666
$ws.Range("B39").Value = 11.41
$ws.Range("C39").Value = 0
$ws.Range("A40").Value = 45919.39583333334
$ws.Range("B40").Value = 4.369
$ws.Range("C40").Value = 3.424
$ws.Range("A41").Value = 45919.40625
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 26.476
$ws.Range("A42").Value = 45919.41666666666
$ws.Range("B42").Value = 0.296
$ws.Range("C42").Value = 6.956
$ws.Rows("43:43").Delete()

Write-Host "done"
